$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L3").Value = 22.22
$ws.Range("Q3").Value = 50.53
$ws.Range("S3").Value = 33.33

$ws.Range("G4").Value = 70
$ws.Range("L4").Value = 25
$ws.Range("Q4").Value = 56
$ws.Range("S4").Value = 27.27
$ws.Range("AE4").Value = 92.67
